$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Thigh Distance [cm] (B4) and Shank Distance [cm] (B8)
$ws.Range("B4").Value = 53
$ws.Range("B8").Value = 27

# Update the active cell selection to B5
$ws.Range("B5").Select()
